$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.00729022,
    0.0042787800000000003,
    0.01412152,
    0.0026985400000000001,
    0.0035630000000000002,
    0.013618079999999999,
    0.0087738999999999994,
    0.0080981400000000002,
    0.014612119999999999,
    0.10361328,
    0.050498439999999999,
    0.076046520000000006,
    0.0022537400000000002,
    0.099354559999999995,
    0.0035768599999999998,
    1.34262456,
    3.3065276799999999,
    0.016119580000000001
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("A3:B20").Select() | Out-Null
